$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '30.768.87'
$ws.Range("E2").Value = '  +0.64%  '

# Row 3
$ws.Range("D3").Value = '1.961.98'
$ws.Range("E3").Value = '  +1.99%  '

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.9998'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  -0.03%  '

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '248.72'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +1.37%  '

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '1.002'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +0.28%  '

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4834'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  -0.75%  '

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '44.63'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +1.64%  '

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.2932'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +0.92%  '

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.06758'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +0.48%  '

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '109.09'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -1.57%  '

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '19.05'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -0.39%  '

# Row 13
$ws.Range("D13").Value = '1.961.04'
$ws.Range("E13").Value = '  +2.03%  '

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.07740'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +2.11%  '

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '5.462'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +3.08%  '

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.6977'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +4.09%  '

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '292.33'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -1.34%  '

# Row 18
$ws.Range("D18").Value = '30.800.06'
$ws.Range("E18").Value = '  +0.81%  '

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '5.663'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +1.87%  '

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '13.15'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +0.90%  '

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '0.000007698'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +1.56%  '

# Row 22
$ws.Range("D22").Value = '2.220.99'
$ws.Range("E22").Value = '  +2.56%  '

# Row 23
$ws.Range("B23").Value = 'BitDAO'
$ws.Range("C23").Value = 'https://coinranking.com/coin/N2IgQ9Xme+bitdao-bit'
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.4916'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +12.24%  '

# Row 24
$ws.Range("B24").Value = 'Dai'
$ws.Range("C24").Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range("E24").Value = '  +0.30%  '

# Row 25
$ws.Range("B25").Value = 'BinanceUSD'
$ws.Range("C25").Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.9989'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -0.16%  '

# Row 26
$ws.Range("B26").Value = 'Chainlink'
$ws.Range("C26").Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '6.607'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +2.19%  '

# Row 27
$ws.Range("B27").Value = 'Cosmos'
$ws.Range("C27").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '9.893'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +4.48%  '

# Row 28
$ws.Range("B28").Value = 'Monero'
$ws.Range("C28").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '170.30'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +3.41%  '

# Row 29
$ws.Range("B29").Value = 'EthereumClassic'
$ws.Range("C29").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '20.06'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -1.01%  '

# Row 30
$ws.Range("B30").Value = 'LidoDAOToken'
$ws.Range("C30").Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '2.174'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +3.03%  '

# Row 31
$ws.Range("B31").Value = 'Stellar'
$ws.Range("C31").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.1074'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +0.04%  '

# Row 32
$ws.Range("B32").Value = 'Toncoin'
$ws.Range("C32").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '1.447'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -0.30%  '

# Row 33
$ws.Range("B33").Value = 'Filecoin'
$ws.Range("C33").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '4.872'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +19.90%  '

# Row 34
$ws.Range("B34").Value = 'InternetComputer(DFINITY)'
$ws.Range("C34").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '4.478'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +7.70%  '

# Row 35
$ws.Range("B35").Value = 'Hedera'
$ws.Range("C35").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.05106'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +1.43%  '

# Row 36
$ws.Range("B36").Value = 'ImmutableX'
$ws.Range("C36").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.7710'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +4.09%  '

# Row 37
$ws.Range("B37").Value = 'ARBITRUM'
$ws.Range("C37").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '1.179'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +3.39%  '

# Row 38
$ws.Range("B38").Value = 'VeChain'
$ws.Range("C38").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.02046'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +0.94%  '

# Row 39
$ws.Range("B39").Value = 'HuobiToken'
$ws.Range("C39").Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '2.729'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +0.62%  '

# Row 40
$ws.Range("B40").Value = 'MXToken'
$ws.Range("C40").Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '2.721'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +1.34%  '

# Row 41
$ws.Range("B41").Value = 'FraxShare'
$ws.Range("C41").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '6.513'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +11.62%  '

# Row 42
$ws.Range("B42").Value = 'RenderToken'
$ws.Range("C42").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '2.131'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +5.35%  '

# Row 43
$ws.Range("B43").Value = 'TrustWalletToken'
$ws.Range("C43").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.8901'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +2.55%  '

# Row 44
$ws.Range("B44").Value = 'Quant'
$ws.Range("C44").Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '110.11'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -0.35%  '

# Row 45
$ws.Range("B45").Value = 'TheSandbox'
$ws.Range("C45").Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.4463'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +0.58%  '

# Row 46
$ws.Range("E46").Value = '  +0.26%  '

# Row 47
$ws.Range("B47").Value = 'Aave'
$ws.Range("C47").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '69.94'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -1.48%  '

# Row 48
$ws.Range("B48").Value = 'Aptos'
$ws.Range("C48").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '7.464'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +2.66%  '

# Row 49
$ws.Range("B49").Value = 'Algorand'
$ws.Range("C49").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.1275'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +3.35%  '

# Row 50
$ws.Range("B50").Value = 'EnergySwap'
$ws.Range("C50").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '9.366'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +1.18%  '

# Row 51
$ws.Range("B51").Value = 'Elrond'
$ws.Range("C51").Value = 'https://coinranking.com/coin/omwkOTglq+elrond-egld'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '36.02'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +3.20%  '

Write-Output "Applied 154 cell updates"